$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "TP19_EE_E9" variable, right after "TP19_EE_E6" (row 18)
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "TP19_EE_E9"

# Match formatting of the other variable rows (border/alignment style)
$ws.Range("A20").Copy()
$ws.Range("A19").PasteSpecial(-4122)

# Update the current selection / scroll position left by the editor
$excel.Goto($ws.Range("A16"), $true)
[void]$ws.Range("B22").Select()
